# Add a new "Italy" worksheet, cloned from "Slovakia", with its own
# market/ticket values, positioned as the last (7th) sheet and made the
# active tab - mirrors how the other per-country sheets (Germany, Belgium,
# Czech, Swiss, Portugal, Slovakia) were originally produced.

$wb = $excel.ActiveWorkbook

$slovakia = $wb.Worksheets.Item("Slovakia")

# Duplicate the Slovakia sheet and place the copy right after it.
$slovakia.Copy([System.Reflection.Missing]::Value, $slovakia) | Out-Null

# The freshly-created copy is now the last sheet in the workbook.
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Fill in the country-specific data (write the ticket id first so the new
# shared-string entries land in the same order as the source edit: ticket
# id before market name).
$italy.Range("B4").Value = "NGC-3145/T2237"
$italy.Range("B2").Value = "Italy Market"

# Restore Slovakia's on-screen selection to a "whole sheet" selection (it
# loses tab focus now that Italy is the active sheet).
$slovakia.Cells.Select() | Out-Null

# Make Italy the active sheet/tab, with A7 selected.
$italy.Activate() | Out-Null
$italy.Range("A7").Select() | Out-Null
